$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8-74 down to 9-75
$ws.Rows.Item(8).Insert()

# Set content for the newly inserted row 8
$ws.Range("A8").Value = "5 - HEAD OF SPORTS COMMITTEE"
$ws.Range("B8").Value = 31

# Renumber the prefix number of each shifted row (9-75) by incrementing by 1
$ws.Range("A9").Value = "6 - OFFICIAL 1"
$ws.Range("A10").Value = "7 - KSENIA"
$ws.Range("A11").Value = "8 - OFFICIAL 2"
$ws.Range("A12").Value = "9 - PAULAUSKAS"
$ws.Range("A13").Value = "10 - GOMELKIY"
$ws.Range("A14").Value = "11 - ALEXANDER BELOV"
$ws.Range("A15").Value = "12 - KORKIYA"
$ws.Range("A16").Value = "13 - EDESHKO"
$ws.Range("A17").Value = "14 - SAKANDELIDZE"
$ws.Range("A18").Value = "15 - ZHARMUKHAMEDOV"
$ws.Range("A19").Value = "16 - SERGEY BELOV"
$ws.Range("A20").Value = "17 - PLAYER"
$ws.Range("A21").Value = "18 - DOCTOR"
$ws.Range("A22").Value = "19 - PLAYERS"
$ws.Range("A23").Value = "20 - TRANSLATOR"
$ws.Range("A24").Value = "21 - TEAM"
$ws.Range("A25").Value = "22 - PASSAGE"
$ws.Range("A26").Value = "23 - AIRPORT WORKER 1"
$ws.Range("A27").Value = "24 - CHAYKIN"
$ws.Range("A28").Value = "25 - CUSTOMER-FAN"
$ws.Range("A29").Value = "26 - CUSTOMER 1"
$ws.Range("A30").Value = "27 - CUSTOMER 2"
$ws.Range("A31").Value = "28 - KSENIYA"
$ws.Range("A32").Value = "29 - TRAINER"
$ws.Range("A33").Value = "30 - POLYAKOVA"
$ws.Range("A34").Value = "31 - SURKOVA"
$ws.Range("A35").Value = "32 - CANARIS"
$ws.Range("A36").Value = "33 - SMIRNOVA"
$ws.Range("A37").Value = "34 - NIKOLAYEVA"
$ws.Range("A38").Value = "35 - SVESHNIKOVA"
$ws.Range("A39").Value = "36 - VLASOVA"
$ws.Range("A40").Value = "37 - SPORTSWOMAN"
$ws.Range("A41").Value = "38 - SASHA BELOV"
$ws.Range("A42").Value = "39 - EDESHKO/ZHARMUKHAMEDOV/BELOV"
$ws.Range("A43").Value = "40 - VOLNOV"
$ws.Range("A44").Value = "41 - TAXI DRIVER"
$ws.Range("A45").Value = "42 - UNCLE ILIKO"
$ws.Range("A46").Value = "43 - VOLNOV’S WIFE"
$ws.Range("A47").Value = "44 - EDESHKO’S GIRLFRIEND"
$ws.Range("A48").Value = "45 - WOMAN"
$ws.Range("A49").Value = "46 - REPORTER"
$ws.Range("A50").Value = "47 - VSEVOLOD"
$ws.Range("A51").Value = "48 - LITHUANIAN REPORTER"
$ws.Range("A52").Value = "49 - VLEVOLOD"
$ws.Range("A53").Value = "50 - CUBAN"
$ws.Range("A54").Value = "51 - LITHUANIAN DRIVER"
$ws.Range("A55").Value = "52 - DRIVER"
$ws.Range("A56").Value = "53 - LITHUANIAN TOURIST"
$ws.Range("A57").Value = "54 - BRUNDAGE’S VOICE"
$ws.Range("A58").Value = "55 - MEDIC"
$ws.Range("A59").Value = "56 - INTERPRETER"
$ws.Range("A60").Value = "57 - ADMINISTRATOR"
$ws.Range("A61").Value = "58 - ANNOUNCER"
$ws.Range("A62").Value = "59 - SERGEY PAVLOVICH"
$ws.Range("A63").Value = "60 - REFEREE"
$ws.Range("A64").Value = "61 - AKSAKAL"
$ws.Range("A65").Value = "62 - SPECTATORS"
$ws.Range("A66").Value = "63 - GRANDFATHER"
$ws.Range("A67").Value = "64 - BASKETBALL PLAYER"
$ws.Range("A68").Value = "65 - GALYA"
$ws.Range("A69").Value = "66 - YERYOMINA"
$ws.Range("A70").Value = "67 - MALE VOICE"
$ws.Range("A71").Value = "68 - AUNT KORKIYA"
$ws.Range("A72").Value = "69 - FANS"
$ws.Range("A73").Value = "70 - NOT IDENTIFIED"
$ws.Range("A74").Value = "71 - ANATOLIY"
$ws.Range("A75").Value = "72 - YEREMINA"
